$wb = $excel.ActiveWorkbook

# --- Update version label on Instructions sheet (v17.b -> v17.c) ---
$wsInstr = $wb.Worksheets.Item("Instructions")
$wsInstr.Range("B1").Value = "v17.c"

# --- Update Severity-Mortality data table ---
$ws = $wb.Worksheets.Item("Severity-Mortality")

# New age-stratified severity (col B) / mortality (col C) percentages,
# rows 2..22 (age categories 0-5y up to 100+y)
$data = @(
    @(0.1, 0.6),
    @(0.1, 0.6),
    @(0.1, 0.6),
    @(0.1, 0.6),
    @(0.5, 1.1),
    @(0.5, 1.1),
    @(1.1, 1.9),
    @(1.1, 1.9),
    @(1.4, 3.3),
    @(1.4, 3.3),
    @(2.9, 6.5),
    @(2.9, 6.5),
    @(5.8, 12.6),
    @(5.8, 12.6),
    @(9.3, 21),
    @(9.3, 21),
    @(26.2, 31.6),
    @(26.2, 31.6),
    @(26.2, 31.6),
    @(26.2, 31.6),
    @(26.2, 31.6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
}

# --- Resize column C to fit the new (narrower) data ---
$ws.Columns.Item(3).ColumnWidth = 14.14

# --- Move the sheet's remembered selection to C22, then restore the
#     Instructions sheet as the active tab (matches the saved workbook) ---
$ws.Range("C22").Select() | Out-Null
$wsInstr.Activate() | Out-Null
